$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Function Registers")

# ---------------------------------------------------------------
# Clear the old Q2:R10 lookup table (values removed, formatting kept)
# ---------------------------------------------------------------
$ws2.Range("Q2:R10").ClearContents()

# ---------------------------------------------------------------
# Servo section (0x20) - registers reworked:
#   "Set Angle"/"Set Speed"/"Get Angle?"/"Get Speed?" (4 regs)
#   -> "Angle" (R/W) / "Busy" (R) (2 regs)
# ---------------------------------------------------------------
$ws2.Range("E20").Value = "Angle"
$ws2.Range("G20").Value = "R/W"
$ws2.Range("K20").Value = "Angle (7::0)"
$ws2.Range("L20").ClearContents()

$ws2.Range("E21").Value = "Busy"
$ws2.Range("F21").Value = "2 bytes"
$ws2.Range("G21").Value = "R"
$ws2.Range("K21").Value = "True/False"
$ws2.Range("L21").ClearContents()

$ws2.Range("C22:J22").ClearContents()
$ws2.Range("C23:J23").ClearContents()

# ---------------------------------------------------------------
# Relay section (0x22): "Value" register renamed to "Bulk",
# plus a new "Channel" (0x02) register added in row 31
# ---------------------------------------------------------------
$ws2.Range("E30").Value = "Bulk"
$ws2.Range("J30").Value = "Value"
$ws2.Range("K30").ClearContents()

$ws2.Range("C30:N30").Copy()
$ws2.Range("C31:N31").PasteSpecial(-4122)
$ws2.Range("C31").Value = "0x22"
$ws2.Range("D31").Value = "0x02"
$ws2.Range("E31").Value = "Channel"
$ws2.Range("F31").Value = "3 Bytes"
$ws2.Range("G31").Value = "R/W"
$ws2.Range("J31").Value = "Channel"
$ws2.Range("K31").Value = "Value"

# ---------------------------------------------------------------
# IR section (0x23): existing Send/Receive rows get Base-Register
# columns (H/I) filled in, plus two new registers appended
# (Send Status / Receive Status)
# ---------------------------------------------------------------
$ws2.Range("H34").Value = "0x23"
$ws2.Range("I34").Value = "0x01"

$ws2.Range("E35").Value = "Send Status"
$ws2.Range("H35").Value = "0x23"
$ws2.Range("I35").Value = "0x02"
$ws2.Range("K35").Value = "Pass/Fail"

$ws2.Range("C35:N35").Copy()
$ws2.Range("C36:N37").PasteSpecial(-4122)

$ws2.Range("C36").Value = "0x23"
$ws2.Range("D36").Value = "0x03"
$ws2.Range("E36").Value = "Receive Status"
$ws2.Range("F36").Value = "2 Bytes"
$ws2.Range("G36").Value = "R"
$ws2.Range("H36").Value = "0x23"
$ws2.Range("I36").Value = "0x03"
$ws2.Range("J36").Value = "Address"
$ws2.Range("K36").Value = "Buffer Size"

$ws2.Range("C37").Value = "0x23"
$ws2.Range("D37").Value = "0x04"
$ws2.Range("E37").Value = "Receive"
$ws2.Range("F37").Value = "2 Bytes"
$ws2.Range("G37").Value = "R"
$ws2.Range("H37").Value = "0x23"
$ws2.Range("I37").Value = "0x04"
$ws2.Range("J37").Value = "Address"
$ws2.Range("K37").Value = "Data"

# ---------------------------------------------------------------
# New "Current - 0x30" register block (client registration screen)
# ---------------------------------------------------------------
$ws2.Range("C29:N29").Copy()
$ws2.Range("C42:N42").PasteSpecial(-4122)
$ws2.Range("D42").Value = "Current - 0x30"
$ws2.Range("J42").Value = 2
$ws2.Range("K42").Value = 3
$ws2.Range("L42").Value = 4
$ws2.Range("M42").Value = 5
$ws2.Range("N42").Value = 6
$ws2.Range("D42:G42").Merge()

$ws2.Range("C30:N30").Copy()
$ws2.Range("C43:N45").PasteSpecial(-4122)

$ws2.Range("C43").Value = "0x30"
$ws2.Range("D43").Value = "0x01"
$ws2.Range("E43").Value = "Value"
$ws2.Range("F43").Value = "2 Bytes"
$ws2.Range("G43").Value = "R/W"
$ws2.Range("J43").Value = "Channel"
$ws2.Range("K43").Value = "Value"

$ws2.Range("C44").Value = "0x30"
$ws2.Range("D44").Value = "0x02"
$ws2.Range("E44").Value = "Trigger Low"
$ws2.Range("F44").Value = "2 Bytes"
$ws2.Range("G44").Value = "R/W"
$ws2.Range("J44").Value = "Channel"
$ws2.Range("K44").Value = "Value"

$ws2.Range("C45").Value = "0x30"
$ws2.Range("D45").Value = "0x03"
$ws2.Range("E45").Value = "Trigger High"
$ws2.Range("F45").Value = "2 Bytes"
$ws2.Range("G45").Value = "R/W"
$ws2.Range("J45").Value = "Channel"
$ws2.Range("K45").Value = "Value"

# ---------------------------------------------------------------
# Selection moves to the newly added row
# ---------------------------------------------------------------
$ws2.Range("D31").Select()
